$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 937.14703
$ws.Range("I15").Value = 937.14703
$ws.Range("K15").Value = 2811.44109
$ws.Range("M15").Value = -2642.44109
$ws.Range("H131").Value = 1874.9231
$ws.Range("I131").Value = 692.9286
$ws.Range("J131").Value = 3253.9167
$ws.Range("K131").Value = 2078.7858
$ws.Range("L131").Value = 9761.750100000001
$ws.Range("M131").Value = 2961.2142
$ws.Range("N131").Value = -19841.7501
$ws.Range("H135").Value = 438.25
$ws.Range("I135").Value = 422.21054
$ws.Range("K135").Value = 3799.89486
$ws.Range("M135").Value = -1264.89486
$ws.Range("H137").Value = 28895.473
$ws.Range("I137").Value = 845
$ws.Range("K137").Value = 2535
$ws.Range("M137").Value = 15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 347879.88
$ws.Range("I2").Value = 427988.7
$ws.Range("J2").Value = 741.6667
$ws.Range("K2").Value = 427988.7
$ws.Range("L2").Value = 741.6667
$ws.Range("M2").Value = -427875.7
$ws.Range("N2").Value = -967.6667
$ws.Range("H32").Value = 3151.0618
$ws.Range("I32").Value = 2604.1023
$ws.Range("J32").Value = 8499.111000000001
$ws.Range("K32").Value = 2604.1023
$ws.Range("L32").Value = 8499.111000000001
$ws.Range("M32").Value = -2317.1023
$ws.Range("N32").Value = -9073.111000000001
$ws.Range("H61").Value = 3871.4783
$ws.Range("I61").Value = 1169.1666
$ws.Range("J61").Value = 13599.8
$ws.Range("K61").Value = 1169.1666
$ws.Range("L61").Value = 13599.8
$ws.Range("M61").Value = -957.1666
$ws.Range("N61").Value = -14023.8
$ws.Range("H74").Value = 2040.0667
$ws.Range("I74").Value = 1515.1904
$ws.Range("J74").Value = 3264.7778
$ws.Range("K74").Value = 1515.1904
$ws.Range("L74").Value = 3264.7778
$ws.Range("M74").Value = -641.1904
$ws.Range("N74").Value = -5012.7778
$ws.Range("H77").Value = 2040.0667
$ws.Range("I77").Value = 1515.1904
$ws.Range("J77").Value = 3264.7778
$ws.Range("K77").Value = 7575.951999999999
$ws.Range("L77").Value = 16323.889
$ws.Range("M77").Value = -3207.951999999999
$ws.Range("N77").Value = -25059.889
$ws.Range("H110").Value = 197.9
$ws.Range("I110").Value = 133
$ws.Range("J110").Value = 349.33334
$ws.Range("K110").Value = 133
$ws.Range("L110").Value = 349.33334
$ws.Range("M110").Value = 1912
$ws.Range("N110").Value = -4439.33334
$ws.Range("H116").Value = 347879.88
$ws.Range("I116").Value = 427988.7
$ws.Range("J116").Value = 741.6667
$ws.Range("K116").Value = 427988.7
$ws.Range("L116").Value = 741.6667
$ws.Range("M116").Value = -425694.7
$ws.Range("N116").Value = -5329.6667
$ws.Range("H132").Value = 1383.9318
$ws.Range("I132").Value = 1251.1143
$ws.Range("J132").Value = 1900.4445
$ws.Range("K132").Value = 3753.3429
$ws.Range("L132").Value = 5701.333500000001
$ws.Range("M132").Value = -1223.3429
$ws.Range("N132").Value = -10761.3335
$ws.Range("H136").Value = 3871.4783
$ws.Range("I136").Value = 1169.1666
$ws.Range("J136").Value = 13599.8
$ws.Range("K136").Value = 3507.4998
$ws.Range("L136").Value = 40799.39999999999
$ws.Range("M136").Value = -957.4998000000001
$ws.Range("N136").Value = -45899.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 347879.88
$ws.Range("I3").Value = 427988.7
$ws.Range("J3").Value = 741.6667
$ws.Range("K3").Value = 427988.7
$ws.Range("L3").Value = 741.6667
$ws.Range("M3").Value = -427874.7
$ws.Range("N3").Value = -969.6667
$ws.Range("H105").Value = 2359.348
$ws.Range("I105").Value = 2343.25
$ws.Range("J105").Value = 2466.6667
$ws.Range("K105").Value = 2343.25
$ws.Range("L105").Value = 2466.6667
$ws.Range("M105").Value = -596.25
$ws.Range("N105").Value = -5960.6667
$ws.Range("H134").Value = 4745.5757
$ws.Range("I134").Value = 5013.0356
$ws.Range("J134").Value = 3247.8
$ws.Range("K134").Value = 15039.1068
$ws.Range("L134").Value = 9743.400000000001
$ws.Range("M134").Value = -12504.1068
$ws.Range("N134").Value = -14813.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2388.9285
$ws.Range("I31").Value = 1959.4445
$ws.Range("J31").Value = 3162
$ws.Range("K31").Value = 1959.4445
$ws.Range("L31").Value = 3162
$ws.Range("M31").Value = -1664.4445
$ws.Range("N31").Value = -3752
$ws.Range("H34").Value = 2388.9285
$ws.Range("I34").Value = 1959.4445
$ws.Range("J34").Value = 3162
$ws.Range("K34").Value = 1959.4445
$ws.Range("L34").Value = 3162
$ws.Range("M34").Value = -1757.4445
$ws.Range("N34").Value = -3566
$ws.Range("H58").Value = 1012322.25
$ws.Range("I58").Value = 1553711.5
$ws.Range("J58").Value = 1729.0667
$ws.Range("K58").Value = 1553711.5
$ws.Range("L58").Value = 1729.0667
$ws.Range("M58").Value = -1553508.5
$ws.Range("N58").Value = -2135.0667
$ws.Range("H96").Value = 31292.25
$ws.Range("J96").Value = 31292.25
$ws.Range("L96").Value = 31292.25
$ws.Range("N96").Value = -36784.25
$ws.Range("H132").Value = 2024.2162
$ws.Range("I132").Value = 1325.2
$ws.Range("J132").Value = 5020
$ws.Range("K132").Value = 3975.6
$ws.Range("L132").Value = 15060
$ws.Range("M132").Value = -1445.6
$ws.Range("N132").Value = -20120
$ws.Range("H134").Value = 1782.0408
$ws.Range("I134").Value = 1675.8918
$ws.Range("J134").Value = 2109.3333
$ws.Range("K134").Value = 5027.6754
$ws.Range("L134").Value = 6327.999899999999
$ws.Range("M134").Value = -2492.6754
$ws.Range("N134").Value = -11397.9999
$ws.Range("H136").Value = 1012322.25
$ws.Range("I136").Value = 1553711.5
$ws.Range("J136").Value = 1729.0667
$ws.Range("K136").Value = 4661134.5
$ws.Range("L136").Value = 5187.2001
$ws.Range("M136").Value = -4658584.5
$ws.Range("N136").Value = -10287.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7138.5
$ws.Range("I56").Value = 7138.5
$ws.Range("K56").Value = 7138.5
$ws.Range("M56").Value = -6608.5
$ws.Range("H107").Value = 453.1905
$ws.Range("I107").Value = 360.5
$ws.Range("J107").Value = 510.23077
$ws.Range("K107").Value = 1081.5
$ws.Range("L107").Value = 1530.69231
$ws.Range("M107").Value = 838.5
$ws.Range("N107").Value = -5370.69231
$ws.Range("H113").Value = 40064.855
$ws.Range("I113").Value = 275476
$ws.Range("J113").Value = 829.6667
$ws.Range("K113").Value = 826428
$ws.Range("L113").Value = 2489.0001
$ws.Range("M113").Value = -824258
$ws.Range("N113").Value = -6829.0001
$ws.Range("H118").Value = 1674.125
$ws.Range("I118").Value = 800
$ws.Range("J118").Value = 1965.5
$ws.Range("K118").Value = 2400
$ws.Range("L118").Value = 5896.5
$ws.Range("M118").Value = -1157
$ws.Range("N118").Value = -8382.5
$ws.Range("H132").Value = 1225.9546
$ws.Range("I132").Value = 781.3333
$ws.Range("J132").Value = 1392.6875
$ws.Range("K132").Value = 7031.9997
$ws.Range("L132").Value = 12534.1875
$ws.Range("M132").Value = -4501.9997
$ws.Range("N132").Value = -17594.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3401.2727
$ws.Range("I70").Value = 3157.3333
$ws.Range("K70").Value = 3157.3333
$ws.Range("M70").Value = -2887.3333
$ws.Range("H73").Value = 3401.2727
$ws.Range("I73").Value = 3157.3333
$ws.Range("K73").Value = 3157.3333
$ws.Range("M73").Value = -2221.3333
$ws.Range("H92").Value = 21127.428
$ws.Range("J92").Value = 21127.428
$ws.Range("L92").Value = 21127.428
$ws.Range("N92").Value = -24871.428
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H132").Value = 593366.0600000001
$ws.Range("I132").Value = 875313.25
$ws.Range("K132").Value = 2625939.75
$ws.Range("M132").Value = -2623409.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 12000
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12386
$ws.Range("H104").Value = 20592.5
$ws.Range("J104").Value = 20592.5
$ws.Range("L104").Value = 20592.5
$ws.Range("N104").Value = -27580.5
$ws.Range("H106").Value = 20400
$ws.Range("J106").Value = 20400
$ws.Range("L106").Value = 20400
$ws.Range("N106").Value = -22924
$ws.Range("H132").Value = 1534.4318
$ws.Range("I132").Value = 1056.5084
$ws.Range("J132").Value = 2506.7585
$ws.Range("K132").Value = 3169.5252
$ws.Range("L132").Value = 7520.2755
$ws.Range("M132").Value = -639.5252
$ws.Range("N132").Value = -12580.2755
$ws.Range("H136").Value = 2453.7844
$ws.Range("I136").Value = 1813.3096
$ws.Range("J136").Value = 5442.6665
$ws.Range("K136").Value = 5439.9288
$ws.Range("L136").Value = 16327.9995
$ws.Range("M136").Value = -2889.9288
$ws.Range("N136").Value = -21427.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 35700
$ws.Range("J69").Value = 35700
$ws.Range("L69").Value = 35700
$ws.Range("N69").Value = -37198
$ws.Range("H72").Value = 35700
$ws.Range("J72").Value = 35700
$ws.Range("L72").Value = 107100
$ws.Range("N72").Value = -114588
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
$ws.Range("H132").Value = 1169.1731
$ws.Range("I132").Value = 689.675
$ws.Range("K132").Value = 2069.025
$ws.Range("M132").Value = 460.9750000000004
$ws.Range("H135").Value = 102907.375
$ws.Range("J135").Value = 102907.375
$ws.Range("L135").Value = 102907.375
$ws.Range("N135").Value = -113047.375
$ws.Range("H136").Value = 15017093
$ws.Range("I136").Value = 16837108
$ws.Range("K136").Value = 50511324
$ws.Range("M136").Value = -50508774
